$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The report gained one more low-stock item: "ZADOGLOBIN 20 CAPS".
# It sorts alphabetically right before "ZYRTEC 10MG/ML ORAL DROPS 10 ML",
# which currently lives in row 18. So: insert a blank row at 18 (this
# pushes the old row 18..24 content down to 19..25, and fixes up the
# merged cells automatically), clone the formatting of the row that is
# now directly below (row 19, a twin of the old row 18 template), fill
# in the new item's data, bump the grand total, and refresh the
# "printed at" timestamp in the footer.
# ------------------------------------------------------------------

# 1) Insert a new blank row at position 18, shifting rows 18-24 down to 19-25.
$ws.Rows("18:18").Insert()

# 2) Clone formatting (styles + merges) for the new row from row 19,
#    which is built from the same per-column style template as every
#    other data row in the table.
$ws.Range("A19:Q19").Copy($ws.Range("A18:Q18"))

# 3) Populate the new row 18 with the ZADOGLOBIN data (item #12 in the
#    renumbered list). Values that look numeric but are stored as plain
#    text in this report use a leading apostrophe so Excel keeps them
#    as text instead of silently re-typing them as numbers.
$ws.Range("A18").Value = 12
$ws.Range("C18").Value = "ZADOGLOBIN 20 CAPS"
$ws.Range("H18").Value = "'1:0"
$ws.Range("L18").Value = "'1"
$ws.Range("N18").Value = "'135.00"
$ws.Range("P18").Value = "'67.5000"
$ws.Range("Q18").Value = "'0:1"

# 4) Row heights: row 18 is a brand-new row (use the same height it had
#    originally at that slot); rows 19-25 keep the heights that already
#    belong to those row slots after the shift, except the grand-total
#    row (now 24) and the footer row (now 25) which need to be set
#    explicitly to match the regenerated report.
$ws.Rows("18:18").RowHeight = 24.75
$ws.Rows("19:19").RowHeight = 25.5
$ws.Rows("20:20").RowHeight = 24.75
$ws.Rows("21:21").RowHeight = 25.5
$ws.Rows("22:22").RowHeight = 25.5
$ws.Rows("23:23").RowHeight = 24.75
$ws.Rows("24:24").RowHeight = 25.5
$ws.Rows("25:25").RowHeight = 16.5

# 5) Update the grand total (old total 415.595 + the new item's 67.5000).
$ws.Range("P24").Value = 483.095

# 6) Refresh the "printed at" timestamp in the footer.
$ws.Range("A25").Value = "Thursday, 21 August, 2025 11:16 AM"
